$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.918.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.711.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.67%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("E6").Value = '  +8.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '656.91'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("E9").Value = '  +4.34%  '
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.710.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.53%  '
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.401.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000271'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.916.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.713.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.530'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '527.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.11%  '
$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("E25").Value = '  +3.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000206'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.31%  '
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.32%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +15.53%  '
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.87%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  +6.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.607'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +19.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.164'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.977'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.17%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0461'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.444'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.86%  '
